$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row number -> new Price (column D) / Volume(1h) (column E) text values.
# Column D values must be written as literal text (prices use '.' both as
# thousands separator and decimal point, and some have significant
# trailing zeros), so we force a text number format before assigning and
# restore the default style afterwards so no new cell style is introduced.
$updates = @{
    2  = @("30.409.38", "  -1.29%  ")
    3  = @("1.904.17", "  -1.83%  ")
    4  = @("0.9997", "  -0.04%  ")
    5  = @("238.53", "  -1.50%  ")
    6  = @("0.9995", "  -0.06%  ")
    7  = @("0.4747", "  -3.11%  ")
    8  = @("0.2828", "  -3.27%  ")
    9  = @("0.06683", "  -2.69%  ")
    10 = @("19.48", "  +1.31%  ")
    11 = @("102.35", "  -2.83%  ")
    12 = @("0.07745", "  -0.10%  ")
    13 = @("1.905.35", "  -1.82%  ")
    14 = @("5.191", "  -2.99%  ")
    15 = @("0.6672", "  -4.22%  ")
    16 = @("277.13", "  +0.81%  ")
    17 = @("30.450.02", "  -1.16%  ")
    18 = @("0.9996", "  +0.01%  ")
    19 = @("0.000007452", "  -3.17%  ")
    20 = @("12.66", "  -2.91%  ")
    21 = @("5.389", "  -2.32%  ")
    22 = @("0.4635", "  -8.33%  ")
    23 = @("0.9999", "  +0.01%  ")
    24 = @("6.256", "  -4.19%  ")
    25 = @("9.307", "  -4.71%  ")
    26 = @("166.85", "  +0.36%  ")
    27 = @("19.20", "  -2.01%  ")
    28 = @("2.059", "  -4.26%  ")
    29 = @("1.379", "  -0.78%  ")
    30 = @("0.09985", "  -3.57%  ")
    31 = @("4.627", "  +1.57%  ")
    32 = @("1.513", "  -2.53%  ")
    33 = @("4.231", "  -2.75%  ")
    34 = @("0.04673", "  -4.14%  ")
    35 = @("0.7302", "  -2.86%  ")
    36 = @("1.109", "  -3.59%  ")
    37 = @("2.710", "  -0.78%  ")
    38 = @("0.01912", "  -4.04%  ")
    39 = @("2.608", "  -2.13%  ")
    40 = @("6.251", "  -4.24%  ")
    41 = @("74.76", "  -4.88%  ")
    42 = @("1.951", "  -5.66%  ")
    43 = @("0.8562", "  -5.49%  ")
    44 = @("105.09", "  -2.29%  ")
    45 = @("0.4239", "  -3.77%  ")
    46 = @("0.9990", "  +0.05%  ")
    47 = @("7.439", "  -2.57%  ")
    48 = @("954.52", "  -2.24%  ")
    49 = @("0.1197", "  -3.28%  ")
    50 = @("34.55", "  -3.69%  ")
    51 = @("8.745", "  -4.80%  ")
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]

    $priceCell = $ws.Cells.Item($row, 4)
    $priceCell.NumberFormat = "@"
    $priceCell.Value = $vals[0]
    $priceCell.Style = "Normal"

    $ws.Cells.Item($row, 5).Value = $vals[1]
}

# Rows 48 and 49 also swap coin identity: the Algorand/Maker rows
# exchanged places in the ranking (name + link), in addition to the
# price/volume changes already applied above.
$ws.Cells.Item(48, 2).Value = "Maker"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"

$ws.Cells.Item(49, 2).Value = "Algorand"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
